$d = $word.ActiveDocument

$d.Content.Find.Execute("2023-12-14 Thursday", $true, $false, $false, $false, $false, $true, 1, $false, "2023-12-15 Friday", 2) | Out-Null
$d.Content.Find.Execute("67×78=", $true, $false, $false, $false, $false, $true, 1, $false, "47×51=", 2) | Out-Null
$d.Content.Find.Execute("98×87=", $true, $false, $false, $false, $false, $true, 1, $false, "31×26=", 2) | Out-Null
$d.Content.Find.Execute("58×73=", $true, $false, $false, $false, $false, $true, 1, $false, "52×24=", 2) | Out-Null
$d.Content.Find.Execute("44×16=", $true, $false, $false, $false, $false, $true, 1, $false, "91×47=", 2) | Out-Null
$d.Content.Find.Execute("84×74=", $true, $false, $false, $false, $false, $true, 1, $false, "17×93=", 2) | Out-Null
$d.Content.Find.Execute("23×99=", $true, $false, $false, $false, $false, $true, 1, $false, "74×11=", 2) | Out-Null
$d.Content.Find.Execute("51×71=", $true, $false, $false, $false, $false, $true, 1, $false, "29×36=", 2) | Out-Null
$d.Content.Find.Execute("77×87=", $true, $false, $false, $false, $false, $true, 1, $false, "50×39=", 2) | Out-Null
$d.Content.Find.Execute("62×12=", $true, $false, $false, $false, $false, $true, 1, $false, "20×94=", 2) | Out-Null
$d.Content.Find.Execute("65×39=", $true, $false, $false, $false, $false, $true, 1, $false, "47×71=", 2) | Out-Null
$d.Content.Find.Execute("81×25=", $true, $false, $false, $false, $false, $true, 1, $false, "15×15=", 2) | Out-Null
$d.Content.Find.Execute("85×70=", $true, $false, $false, $false, $false, $true, 1, $false, "72×35=", 2) | Out-Null
$d.Content.Find.Execute("31×13=", $true, $false, $false, $false, $false, $true, 1, $false, "51×87=", 2) | Out-Null
$d.Content.Find.Execute("99×62=", $true, $false, $false, $false, $false, $true, 1, $false, "60×62=", 2) | Out-Null
$d.Content.Find.Execute("31×68=", $true, $false, $false, $false, $false, $true, 1, $false, "26×20=", 2) | Out-Null
$d.Content.Find.Execute("88×45=", $true, $false, $false, $false, $false, $true, 1, $false, "48×27=", 2) | Out-Null
$d.Content.Find.Execute("88×73=", $true, $false, $false, $false, $false, $true, 1, $false, "67×32=", 2) | Out-Null
$d.Content.Find.Execute("16×57=", $true, $false, $false, $false, $false, $true, 1, $false, "18×32=", 2) | Out-Null
$d.Content.Find.Execute("12×21=", $true, $false, $false, $false, $false, $true, 1, $false, "66×30=", 2) | Out-Null
$d.Content.Find.Execute("23×92=", $true, $false, $false, $false, $false, $true, 1, $false, "59×41=", 2) | Out-Null
$d.Content.Find.Execute("75×19=", $true, $false, $false, $false, $false, $true, 1, $false, "85×83=", 2) | Out-Null
$d.Content.Find.Execute("67×24=", $true, $false, $false, $false, $false, $true, 1, $false, "86×27=", 2) | Out-Null
$d.Content.Find.Execute("95×84=", $true, $false, $false, $false, $false, $true, 1, $false, "70×14=", 2) | Out-Null
$d.Content.Find.Execute("33×51=", $true, $false, $false, $false, $false, $true, 1, $false, "62×44=", 2) | Out-Null
$d.Content.Find.Execute("17×87=", $true, $false, $false, $false, $false, $true, 1, $false, "83×33=", 2) | Out-Null
